$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 40001944
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 45456590
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 45456590
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -45456940

$ws.Range("H62").Value = 55907.95
$ws.Range("I62").Value = 147328.42
$ws.Range("J62").Value = 6681.5386
$ws.Range("K62").Value = 147328.42
$ws.Range("L62").Value = 6681.5386
$ws.Range("M62").Value = -146704.42
$ws.Range("N62").Value = -7929.5386

$ws.Range("H65").Value = 55907.95
$ws.Range("I65").Value = 147328.42
$ws.Range("J65").Value = 6681.5386
$ws.Range("K65").Value = 736642.1000000001
$ws.Range("L65").Value = 33407.693
$ws.Range("M65").Value = -733522.1000000001
$ws.Range("N65").Value = -39647.693

$ws.Range("H92").Value = 945.2941
$ws.Range("I92").Value = 971.3333
$ws.Range("J92").Value = 750
$ws.Range("K92").Value = 971.3333
$ws.Range("L92").Value = 750
$ws.Range("M92").Value = 276.6667
$ws.Range("N92").Value = -3246

$ws.Range("H132").Value = 1268.5306
$ws.Range("I132").Value = 1372
$ws.Range("J132").Value = 738.25
$ws.Range("K132").Value = 4116
$ws.Range("L132").Value = 2214.75
$ws.Range("M132").Value = -1586
$ws.Range("N132").Value = -7274.75

$ws.Range("H137").Value = 926.38464
$ws.Range("I137").Value = 798
$ws.Range("K137").Value = 2394
$ws.Range("M137").Value = 156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1900
$ws.Range("I45").Value = 1280
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 1280
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -903
$ws.Range("N45").Value = -5754

$ws.Range("H61").Value = 1348.5358
$ws.Range("I61").Value = 802.25
$ws.Range("J61").Value = 2714.25
$ws.Range("K61").Value = 802.25
$ws.Range("L61").Value = 2714.25
$ws.Range("M61").Value = -590.25
$ws.Range("N61").Value = -3138.25

$ws.Range("H136").Value = 1348.5358
$ws.Range("I136").Value = 802.25
$ws.Range("J136").Value = 2714.25
$ws.Range("K136").Value = 2406.75
$ws.Range("L136").Value = 8142.75
$ws.Range("M136").Value = 143.25
$ws.Range("N136").Value = -13242.75

$ws.Range("H139").Value = 20500
$ws.Range("J139").Value = 20500
$ws.Range("L139").Value = 20500
$ws.Range("N139").Value = -30780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 45937.24
$ws.Range("I134").Value = 2396.9678
$ws.Range("K134").Value = 7190.903399999999
$ws.Range("M134").Value = -4655.903399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 69343.336
$ws.Range("J23").Value = 69343.336
$ws.Range("L23").Value = 69343.336
$ws.Range("N23").Value = -69823.336

$ws.Range("H27").Value = 69343.336
$ws.Range("J27").Value = 69343.336
$ws.Range("L27").Value = 69343.336
$ws.Range("N27").Value = -69727.336

$ws.Range("H31").Value = 2229.8833
$ws.Range("I31").Value = 2297.476
$ws.Range("J31").Value = 2072.1667
$ws.Range("K31").Value = 2297.476
$ws.Range("L31").Value = 2072.1667
$ws.Range("M31").Value = -2002.476
$ws.Range("N31").Value = -2662.1667

$ws.Range("H34").Value = 2229.8833
$ws.Range("I34").Value = 2297.476
$ws.Range("J34").Value = 2072.1667
$ws.Range("K34").Value = 2297.476
$ws.Range("L34").Value = 2072.1667
$ws.Range("M34").Value = -2095.476
$ws.Range("N34").Value = -2476.1667

$ws.Range("H132").Value = 1046.3096
$ws.Range("I132").Value = 848.57355
$ws.Range("J132").Value = 1886.6875
$ws.Range("K132").Value = 2545.72065
$ws.Range("L132").Value = 5660.0625
$ws.Range("M132").Value = -15.72064999999975
$ws.Range("N132").Value = -10720.0625

$ws.Range("H134").Value = 1114.875
$ws.Range("I134").Value = 904.2
$ws.Range("J134").Value = 1976.7273
$ws.Range("K134").Value = 2712.6
$ws.Range("L134").Value = 5930.1819
$ws.Range("M134").Value = -177.6000000000004
$ws.Range("N134").Value = -11000.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H96").Value = 4740
$ws.Range("J96").Value = 4740
$ws.Range("L96").Value = 14220
$ws.Range("N96").Value = -18338

$ws.Range("H101").Value = 6666.6665
$ws.Range("J101").Value = 6666.6665
$ws.Range("L101").Value = 19999.9995
$ws.Range("N101").Value = -24867.9995

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H103").Value = 2545
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 3293.3333
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 9879.999899999999
$ws.Range("M103").Value = -21
$ws.Range("N103").Value = -11637.9999

$ws.Range("H105").Value = 156142.86
$ws.Range("J105").Value = 156142.86
$ws.Range("L105").Value = 468428.58
$ws.Range("N105").Value = -473670.58

$ws.Range("H107").Value = 648822.9399999999
$ws.Range("I107").Value = 1152.6
$ws.Range("K107").Value = 3457.8
$ws.Range("M107").Value = -1537.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 9000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H46").Value = 12599
$ws.Range("I46").Value = 5998.5
$ws.Range("J46").Value = 25800
$ws.Range("K46").Value = 5998.5
$ws.Range("L46").Value = 25800
$ws.Range("M46").Value = -5842.5
$ws.Range("N46").Value = -26112

$ws.Range("H50").Value = 9000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H57").Value = 14184.5
$ws.Range("I57").Value = 8702.333000000001
$ws.Range("J57").Value = 19666.666
$ws.Range("K57").Value = 8702.333000000001
$ws.Range("L57").Value = 19666.666
$ws.Range("M57").Value = -7882.333000000001
$ws.Range("N57").Value = -21306.666

$ws.Range("H132").Value = 1788.5
$ws.Range("I132").Value = 1611.8518
$ws.Range("K132").Value = 4835.555399999999
$ws.Range("M132").Value = -2305.555399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2068.6296
$ws.Range("I132").Value = 1906.9546
$ws.Range("J132").Value = 2780
$ws.Range("K132").Value = 5720.8638
$ws.Range("L132").Value = 8340
$ws.Range("M132").Value = -3190.8638
$ws.Range("N132").Value = -13400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1231.5778
$ws.Range("I132").Value = 797.53125
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 2392.59375
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = 137.40625
$ws.Range("N132").Value = -11960

$ws.Range("H136").Value = 904.4545000000001
$ws.Range("I136").Value = 406.125
$ws.Range("J136").Value = 2233.3333
$ws.Range("K136").Value = 1218.375
$ws.Range("L136").Value = 6699.999899999999
$ws.Range("M136").Value = 1331.625
$ws.Range("N136").Value = -11799.9999
